$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Guca2a"
$ws.Range("C2").Value = "Gucy2c"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01873166666666667
$ws.Range("H2").Value = 0.056195
$ws.Range("I2").Value = 0.01714827236381741
$ws.Range("J2").Value = 0.0200504588645364
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3282236666666666
$ws.Range("N2").Value = 0.984671
$ws.Range("O2").Value = 0.1589343097904449
$ws.Range("P2").Value = 0.1589343097904449
$ws.Range("Q2").Value = 0.00614817631611111
$ws.Range("R2").Value = 0.055333586845
$ws.Range("S2").Value = 0.002725448832241881
$ws.Range("T2").Value = 0.0031867058406168

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Guca2a"
$ws.Range("C3").Value = "Gucy2c"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01873166666666667
$ws.Range("H3").Value = 0.056195
$ws.Range("I3").Value = 0.01714827236381741
$ws.Range("J3").Value = 0.0200504588645364
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.716416
$ws.Range("N3").Value = 5.149248
$ws.Range("O3").Value = 0.8311326085767012
$ws.Range("P3").Value = 0.8311326085767011
$ws.Range("Q3").Value = 0.03215133237333333
$ws.Range("R3").Value = 0.28936199136
$ws.Range("S3").Value = 0.01425248834232331
$ws.Range("T3").Value = 0.01666459017924198

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Guca2a"
$ws.Range("C4").Value = "Gucy2c"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01873166666666667
$ws.Range("H4").Value = 0.056195
$ws.Range("I4").Value = 0.01714827236381741
$ws.Range("J4").Value = 0.0200504588645364
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02051333333333333
$ws.Range("N4").Value = 0.06154
$ws.Range("O4").Value = 0.009933081632854
$ws.Range("P4").Value = 0.009933081632853998
$ws.Range("Q4").Value = 0.0003842489222222222
$ws.Range("R4").Value = 0.0034582403
$ws.Range("S4").Value = 0.0001703351892522125
$ws.Range("T4").Value = 0.0001991628446776211

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Guca2a"
$ws.Range("C5").Value = "Gucy2c"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4743275
$ws.Range("H5").Value = 0.948655
$ws.Range("I5").Value = 0.4342324313363431
$ws.Range("J5").Value = 0.3384815028763551
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3282236666666666
$ws.Range("N5").Value = 0.984671
$ws.Range("O5").Value = 0.1589343097904449
$ws.Range("P5").Value = 0.1589343097904449
$ws.Range("Q5").Value = 0.1556855112508333
$ws.Range("R5").Value = 0.934113067505
$ws.Range("S5").Value = 0.06901443176306847
$ws.Range("T5").Value = 0.05379632403648599

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Guca2a"
$ws.Range("C6").Value = "Gucy2c"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4743275
$ws.Range("H6").Value = 0.948655
$ws.Range("I6").Value = 0.4342324313363431
$ws.Range("J6").Value = 0.3384815028763551
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.716416
$ws.Range("N6").Value = 5.149248
$ws.Range("O6").Value = 0.8311326085767012
$ws.Range("P6").Value = 0.8311326085767011
$ws.Range("Q6").Value = 0.81414331024
$ws.Range("R6").Value = 4.88485986144
$ws.Range("S6").Value = 0.3609047333851782
$ws.Range("T6").Value = 0.2813230144405872

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Guca2a"
$ws.Range("C7").Value = "Gucy2c"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4743275
$ws.Range("H7").Value = 0.948655
$ws.Range("I7").Value = 0.4342324313363431
$ws.Range("J7").Value = 0.3384815028763551
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02051333333333333
$ws.Range("N7").Value = 0.06154
$ws.Range("O7").Value = 0.009933081632854
$ws.Range("P7").Value = 0.009933081632853998
$ws.Range("Q7").Value = 0.009730038116666666
$ws.Range("R7").Value = 0.0583802287
$ws.Range("S7").Value = 0.004313266188096566
$ws.Range("T7").Value = 0.00336216439928194

$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "Guca2a"
$ws.Range("C8").Value = "Gucy2c"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.4343753333333333
$ws.Range("H8").Value = 1.303126
$ws.Range("I8").Value = 0.3976574352232747
$ws.Range("J8").Value = 0.4649572783754401
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3282236666666666
$ws.Range("N8").Value = 0.984671
$ws.Range("O8").Value = 0.1589343097904449
$ws.Range("P8").Value = 0.1589343097904449
$ws.Range("Q8").Value = 0.1425722646162222
$ws.Range("R8").Value = 1.283150381546
$ws.Range("S8").Value = 0.06320141000024973
$ws.Range("T8").Value = 0.07389766412064433

$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "Guca2a"
$ws.Range("C9").Value = "Gucy2c"
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.4343753333333333
$ws.Range("H9").Value = 1.303126
$ws.Range("I9").Value = 0.3976574352232747
$ws.Range("J9").Value = 0.4649572783754401
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.716416
$ws.Range("N9").Value = 5.149248
$ws.Range("O9").Value = 0.8311326085767012
$ws.Range("P9").Value = 0.8311326085767011
$ws.Range("Q9").Value = 0.7455687721386667
$ws.Range("R9").Value = 6.710118949248
$ws.Range("S9").Value = 0.3305060614570409
$ws.Range("T9").Value = 0.3864411556529029

$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Guca2a"
$ws.Range("C10").Value = "Gucy2c"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4343753333333333
$ws.Range("H10").Value = 1.303126
$ws.Range("I10").Value = 0.3976574352232747
$ws.Range("J10").Value = 0.4649572783754401
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02051333333333333
$ws.Range("N10").Value = 0.06154
$ws.Range("O10").Value = 0.009933081632854
$ws.Range("P10").Value = 0.009933081632853998
$ws.Range("Q10").Value = 0.008910486004444444
$ws.Range("R10").Value = 0.08019437404
$ws.Range("S10").Value = 0.003949963765984139
$ws.Range("T10").Value = 0.004618458601892868

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Guca2a"
$ws.Range("C11").Value = "Gucy2c"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.164901
$ws.Range("H11").Value = 0.494703
$ws.Range("I11").Value = 0.1509618610765648
$ws.Range("J11").Value = 0.1765107598836685
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3282236666666666
$ws.Range("N11").Value = 0.984671
$ws.Range("O11").Value = 0.1589343097904449
$ws.Range("P11").Value = 0.1589343097904449
$ws.Range("Q11").Value = 0.05412441085699999
$ws.Range("R11").Value = 0.487119697713
$ws.Range("S11").Value = 0.02399301919488487
$ws.Range("T11").Value = 0.0280536157926978

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Guca2a"
$ws.Range("C12").Value = "Gucy2c"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.164901
$ws.Range("H12").Value = 0.494703
$ws.Range("I12").Value = 0.1509618610765648
$ws.Range("J12").Value = 0.1765107598836685
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.716416
$ws.Range("N12").Value = 5.149248
$ws.Range("O12").Value = 0.8311326085767012
$ws.Range("P12").Value = 0.8311326085767011
$ws.Range("Q12").Value = 0.283038714816
$ws.Range("R12").Value = 2.547348433344
$ws.Range("S12").Value = 0.1254693253921589
$ws.Range("T12").Value = 0.1467038483039691

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Guca2a"
$ws.Range("C13").Value = "Gucy2c"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.164901
$ws.Range("H13").Value = 0.494703
$ws.Range("I13").Value = 0.1509618610765648
$ws.Range("J13").Value = 0.1765107598836685
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.02051333333333333
$ws.Range("N13").Value = 0.06154
$ws.Range("O13").Value = 0.009933081632854
$ws.Range("P13").Value = 0.009933081632853998
$ws.Range("Q13").Value = 0.003382669179999999
$ws.Range("R13").Value = 0.03044402262
$ws.Range("S13").Value = 0.001499516489521083
$ws.Range("T13").Value = 0.00175329578700157
